$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format first so numeric-looking
# values (e.g. "302.97") are not auto-converted to real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.089.06"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "2.310.27"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "302.97"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "101.87"
$ws.Range("E6").Value = "  +5.93%  "
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +6.33%  "
$ws.Range("D10").Value = "35.92"
$ws.Range("E10").Value = "  +8.29%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").Value = "  +3.70%  "
$ws.Range("D13").Value = "17.92"
$ws.Range("E13").Value = "  +15.33%  "
$ws.Range("D14").Value = "6.94"
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").Value = "2.687.38"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").Value = "2.366.21"
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("E17").Value = "  +4.14%  "
$ws.Range("D18").Value = "43.041.48"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "12.64"
$ws.Range("E19").Value = "  +8.20%  "
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D22").Value = "67.95"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").Value = "237.64"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +12.78%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E28").Value = "  +11.34%  "
$ws.Range("D29").Value = "34.76"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "167.88"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "5.05"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "17.25"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +3.75%  "
$ws.Range("D37").Value = "0.0698"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("E38").Value = "  +3.48%  "
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").Value = "2.002.47"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  -4.36%  "
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +7.94%  "
$ws.Range("D46").Value = "17.68"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").Value = "2.90"
$ws.Range("E47").Value = "  +4.61%  "
$ws.Range("D48").Value = "56.07"
$ws.Range("E48").Value = "  +7.45%  "
$ws.Range("D49").Value = "2.529.37"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("E51").Value = "  +2.51%  "

# Remove the temporary Text format so the cell style matches the original
# (no explicit style index) while keeping the values as text.
$ws.Range("D2:D51").ClearFormats()
